$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.737.16'
$ws.Range("E2").Value = '  -1.70%  '

$ws.Range("D3").Value = '2.361.66'
$ws.Range("E3").Value = '  +0.41%  '

$ws.Range("E4").Value = '  -0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.40'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.61%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.638'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.82%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.625'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.86%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.08'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.56%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0924'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.57%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.47'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.69%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.01'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.33%  '

$ws.Range("E14").Value = '  -0.06%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.52'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.91%  '

$ws.Range("D16").Value = '2.718.45'
$ws.Range("E16").Value = '  +0.37%  '

$ws.Range("D17").Value = '2.361.68'
$ws.Range("E17").Value = '  -2.29%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.09'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +11.74%  '

$ws.Range("D19").Value = '42.642.28'
$ws.Range("E19").Value = '  -1.84%  '

$ws.Range("E20").Value = '  -1.87%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '76.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.27%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.72'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '265.88'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.32'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -9.66%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.11'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.77%  '

$ws.Range("E26").Value = '  +0.20%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.48'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.20%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.97'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.08%  '

$ws.Range("E29").Value = '  -2.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '175.66'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.27%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.10'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.96%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0901'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '35.35'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -10.06%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.06'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.76%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.133'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.12%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.59'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.78%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0359'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.68%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.95'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.81%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.107'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.83%  '

$ws.Range("E40").Value = '  -8.83%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.52'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.94%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.237'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.43%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.21'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.07%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '120.71'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.73%  '

$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.27%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.65'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +22.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.92'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.29%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.54'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.88%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.26'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.97%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.27'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.81%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.100'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.08%  '
